$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix separator typo in contractor names (comma -> period) ---
$ws.Range("E20").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E52").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E118").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"

$ws.Range("E25").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F25").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E50").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F50").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E75").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F75").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"

$ws.Range("E49").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E64").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"

# --- Fix floating point formatting in "Importe" column (H): 
#     remove "." thousands separator, change "," decimal separator to "." ---
$importeRange = $ws.Range("H2:H153")
$importeRange.NumberFormat = "@"

$ws.Cells.Item(2, 8).Value = "14525.00"
$ws.Cells.Item(3, 8).Value = "6679.80"
$ws.Cells.Item(4, 8).Value = "4997.42"
$ws.Cells.Item(5, 8).Value = "1125.00"
$ws.Cells.Item(6, 8).Value = "1060.00"
$ws.Cells.Item(7, 8).Value = "3292.80"
$ws.Cells.Item(8, 8).Value = "104029.39"
$ws.Cells.Item(9, 8).Value = "1562.70"
$ws.Cells.Item(10, 8).Value = "4589.25"
$ws.Cells.Item(11, 8).Value = "4541.51"
$ws.Cells.Item(12, 8).Value = "6360.12"
$ws.Cells.Item(13, 8).Value = "522.00"
$ws.Cells.Item(14, 8).Value = "961.00"
$ws.Cells.Item(15, 8).Value = "8980.59"
$ws.Cells.Item(16, 8).Value = "480.00"
$ws.Cells.Item(17, 8).Value = "488.00"
$ws.Cells.Item(18, 8).Value = "845.76"
$ws.Cells.Item(19, 8).Value = "19834.69"
$ws.Cells.Item(20, 8).Value = "170.00"
$ws.Cells.Item(21, 8).Value = "1521.93"
$ws.Cells.Item(22, 8).Value = "48.52"
$ws.Cells.Item(23, 8).Value = "0.60"
$ws.Cells.Item(24, 8).Value = "105.00"
$ws.Cells.Item(25, 8).Value = "228.00"
$ws.Cells.Item(26, 8).Value = "4853.87"
$ws.Cells.Item(27, 8).Value = "1706.00"
$ws.Cells.Item(28, 8).Value = "687.28"
$ws.Cells.Item(29, 8).Value = "342.50"
$ws.Cells.Item(30, 8).Value = "659.04"
$ws.Cells.Item(31, 8).Value = "50.00"
$ws.Cells.Item(32, 8).Value = "2595.00"
$ws.Cells.Item(33, 8).Value = "155.66"
$ws.Cells.Item(34, 8).Value = "31986.02"
$ws.Cells.Item(35, 8).Value = "1026.60"
$ws.Cells.Item(36, 8).Value = "466.30"
$ws.Cells.Item(37, 8).Value = "12400.00"
$ws.Cells.Item(38, 8).Value = "167.14"
$ws.Cells.Item(39, 8).Value = "36.04"
$ws.Cells.Item(40, 8).Value = "338.94"
$ws.Cells.Item(41, 8).Value = "2241.75"
$ws.Cells.Item(42, 8).Value = "160.00"
$ws.Cells.Item(43, 8).Value = "2400.00"
$ws.Cells.Item(44, 8).Value = "13160.00"
$ws.Cells.Item(45, 8).Value = "540.00"
$ws.Cells.Item(46, 8).Value = "18144.00"
$ws.Cells.Item(47, 8).Value = "128.00"
$ws.Cells.Item(48, 8).Value = "214.50"
$ws.Cells.Item(49, 8).Value = "4659.60"
$ws.Cells.Item(50, 8).Value = "84.79"
$ws.Cells.Item(51, 8).Value = "1803.62"
$ws.Cells.Item(52, 8).Value = "5030.00"
$ws.Cells.Item(53, 8).Value = "85.00"
$ws.Cells.Item(54, 8).Value = "14937.14"
$ws.Cells.Item(55, 8).Value = "4924.70"
$ws.Cells.Item(56, 8).Value = "132870.00"
$ws.Cells.Item(57, 8).Value = "8999.15"
$ws.Cells.Item(58, 8).Value = "995.41"
$ws.Cells.Item(59, 8).Value = "43.43"
$ws.Cells.Item(60, 8).Value = "2943.70"
$ws.Cells.Item(61, 8).Value = "40.69"
$ws.Cells.Item(62, 8).Value = "144.90"
$ws.Cells.Item(63, 8).Value = "23.00"
$ws.Cells.Item(64, 8).Value = "59.00"
$ws.Cells.Item(65, 8).Value = "77.00"
$ws.Cells.Item(66, 8).Value = "34.59"
$ws.Cells.Item(67, 8).Value = "6795.00"
$ws.Cells.Item(68, 8).Value = "920.80"
$ws.Cells.Item(69, 8).Value = "1125.50"
$ws.Cells.Item(70, 8).Value = "12.00"
$ws.Cells.Item(71, 8).Value = "17.97"
$ws.Cells.Item(72, 8).Value = "450.00"
$ws.Cells.Item(73, 8).Value = "1307.00"
$ws.Cells.Item(74, 8).Value = "7500.06"
$ws.Cells.Item(75, 8).Value = "7644.00"
$ws.Cells.Item(76, 8).Value = "180.00"
$ws.Cells.Item(77, 8).Value = "6500.00"
$ws.Cells.Item(78, 8).Value = "8100.00"
$ws.Cells.Item(79, 8).Value = "30000.00"
$ws.Cells.Item(80, 8).Value = "2400.00"
$ws.Cells.Item(81, 8).Value = "2500.00"
$ws.Cells.Item(82, 8).Value = "2785.00"
$ws.Cells.Item(83, 8).Value = "11000.00"
$ws.Cells.Item(84, 8).Value = "1500.00"
$ws.Cells.Item(85, 8).Value = "327.41"
$ws.Cells.Item(86, 8).Value = "84.57"
$ws.Cells.Item(87, 8).Value = "1718.00"
$ws.Cells.Item(88, 8).Value = "103.04"
$ws.Cells.Item(89, 8).Value = "326653.70"
$ws.Cells.Item(90, 8).Value = "5700.00"
$ws.Cells.Item(91, 8).Value = "300.00"
$ws.Cells.Item(92, 8).Value = "750.00"
$ws.Cells.Item(93, 8).Value = "1400.00"
$ws.Cells.Item(94, 8).Value = "2300.00"
$ws.Cells.Item(95, 8).Value = "7393.99"
$ws.Cells.Item(96, 8).Value = "150.00"
$ws.Cells.Item(97, 8).Value = "400.00"
$ws.Cells.Item(98, 8).Value = "3320.00"
$ws.Cells.Item(99, 8).Value = "1250.00"
$ws.Cells.Item(100, 8).Value = "1000.00"
$ws.Cells.Item(101, 8).Value = "850.00"
$ws.Cells.Item(102, 8).Value = "1000.00"
$ws.Cells.Item(103, 8).Value = "4000.00"
$ws.Cells.Item(104, 8).Value = "2652.00"
$ws.Cells.Item(105, 8).Value = "150.00"
$ws.Cells.Item(106, 8).Value = "1150.00"
$ws.Cells.Item(107, 8).Value = "1250.00"
$ws.Cells.Item(108, 8).Value = "150.00"
$ws.Cells.Item(109, 8).Value = "1010.00"
$ws.Cells.Item(110, 8).Value = "3130.00"
$ws.Cells.Item(111, 8).Value = "2870.00"
$ws.Cells.Item(112, 8).Value = "2000.00"
$ws.Cells.Item(113, 8).Value = "400.00"
$ws.Cells.Item(114, 8).Value = "880.00"
$ws.Cells.Item(115, 8).Value = "18335.00"
$ws.Cells.Item(116, 8).Value = "585.00"
$ws.Cells.Item(117, 8).Value = "595.00"
$ws.Cells.Item(118, 8).Value = "285.00"
$ws.Cells.Item(119, 8).Value = "150.00"
$ws.Cells.Item(120, 8).Value = "120.00"
$ws.Cells.Item(121, 8).Value = "10.08"
$ws.Cells.Item(122, 8).Value = "837.32"
$ws.Cells.Item(123, 8).Value = "6395.00"
$ws.Cells.Item(124, 8).Value = "112.32"
$ws.Cells.Item(125, 8).Value = "190.00"
$ws.Cells.Item(126, 8).Value = "290.00"
$ws.Cells.Item(127, 8).Value = "50.00"
$ws.Cells.Item(128, 8).Value = "724.00"
$ws.Cells.Item(129, 8).Value = "44.00"
$ws.Cells.Item(130, 8).Value = "1350.60"
$ws.Cells.Item(131, 8).Value = "7.26"
$ws.Cells.Item(132, 8).Value = "883.00"
$ws.Cells.Item(133, 8).Value = "804.49"
$ws.Cells.Item(134, 8).Value = "120.00"
$ws.Cells.Item(135, 8).Value = "637.80"
$ws.Cells.Item(136, 8).Value = "117364.20"
$ws.Cells.Item(137, 8).Value = "2616.29"
$ws.Cells.Item(138, 8).Value = "836.11"
$ws.Cells.Item(139, 8).Value = "12630.00"
$ws.Cells.Item(140, 8).Value = "377623.20"
$ws.Cells.Item(141, 8).Value = "670.00"
$ws.Cells.Item(142, 8).Value = "201400.00"
$ws.Cells.Item(143, 8).Value = "34000.00"
$ws.Cells.Item(144, 8).Value = "213896.00"
$ws.Cells.Item(145, 8).Value = "20000.00"
$ws.Cells.Item(146, 8).Value = "207261.00"
$ws.Cells.Item(147, 8).Value = "228984.00"
$ws.Cells.Item(148, 8).Value = "207000.00"
$ws.Cells.Item(149, 8).Value = "204524.00"
$ws.Cells.Item(150, 8).Value = "190000.00"
$ws.Cells.Item(151, 8).Value = "6500.00"
$ws.Cells.Item(152, 8).Value = "140.00"
$ws.Cells.Item(153, 8).Value = "52.50"

$importeRange.Style = "Normal"

